$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "291.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-6.70%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.22%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.021"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.64%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07342"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.96%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.02%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.535"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-7.71%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9226"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.77%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.99%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1213"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.27%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1736"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.32%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04311"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.69%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08621"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.53%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1054"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.07%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001267"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.89%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005918"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.25%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.46%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.00%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.678"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.35%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.85%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.20%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.03940"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.92%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.41%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-7.36%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.85%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003727"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02288"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-5.55%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.04973"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.66%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005405"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "147.13%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007679"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.62%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.03%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007349"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.66%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007948"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.40%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3163"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.52%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006350"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.65%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.08%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02039"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-93.11%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.08%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.08%"